# rysboard BOM - "More Work. RPI can now programm STM32 via GPIO"
# Fill in the "Coils" sheet (header + blank template rows, with one
# hyperlink-style URL in F5) and add a new "RGB LED" component row
# (row 10) on the "Diodes and Semiconductors" sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# 1. "Coils" sheet: populate header row + 14 templated rows
# ---------------------------------------------------------------
$coils = $wb.Worksheets.Item("Coils")

# Header row (row 4)
$coils.Range("A4").Value = "No"
$coils.Range("B4").Value = "TYP"
$coils.Range("C4").Value = "Case"
$coils.Range("D4").Value = "Producers Symbol"
$coils.Range("E4").Value = "Opis"
$coils.Range("F4").Value = "Link"
$coils.Range("G4").Value = "Schematic Symbol"
$coils.Range("H4").Value = "No of parts in schematic"

# Numbered rows 5..18 -> values 1..14
for ($i = 0; $i -lt 14; $i++) {
    $row = 5 + $i
    $coils.Cells.Item($row, 1).Value = $i + 1
}

# Row heights that differ from the default, matching the template rows
$coils.Rows.Item(5).RowHeight = 22.5
$coils.Rows.Item(6).RowHeight = 31.5
$coils.Rows.Item(8).RowHeight = 31.5
$coils.Rows.Item(9).RowHeight = 31.5

# The only populated content cell: F5 holds the coil's store link
$coils.Range("F5").Value = "https://www.tme.eu/pl/details/hpi1040-220/dlawiki-smd-mocy/ferrocore/"

# Style the various cells consistent with the other (already filled)
# BOM sheets: "Case" column uses the text-format style, "Producers
# Symbol" cells use the wrapped/bold style, and F5 uses the hyperlink
# style.
$coils.Range("C5:C18").Style = "Normalny"
$coils.Range("C5:C18").NumberFormat = "@"
$coils.Range("D5").Font.Bold = $true
$coils.Range("D5").WrapText = $true
$coils.Range("D5").VerticalAlignment = -4108
$coils.Range("D6").Font.Size = 24
$coils.Range("D6").Font.Bold = $true
$coils.Range("D6").WrapText = $true
$coils.Range("D6").VerticalAlignment = -4108
$coils.Range("D8").Font.Size = 24
$coils.Range("D8").Font.Bold = $true
$coils.Range("D8").WrapText = $true
$coils.Range("D8").VerticalAlignment = -4108
$coils.Range("D9").Font.Size = 24
$coils.Range("D9").Font.Bold = $true
$coils.Range("D9").WrapText = $true
$coils.Range("D9").VerticalAlignment = -4108
$coils.Range("F5").Font.Underline = $true
$coils.Range("F5").Font.ColorIndex = 10
$coils.Range("F8").Font.Underline = $true
$coils.Range("F9").Font.Underline = $true

$coils.Range("F5").Select()

# ---------------------------------------------------------------
# 2. "Diodes and Semiconductors" sheet: add the RGB LED part (row 10)
# ---------------------------------------------------------------
$diodes = $wb.Worksheets.Item("Diodes and Semiconductors")

$diodes.Range("B10").Value = "RGB LED"
$diodes.Range("C10").Value = "~~Ćustom~~~"
$diodes.Range("D10").Value = "KRTBDWLM32.32-T4U6-JW+AAAB-J3+R4S6-5V"
$diodes.Range("F10").Value = "https://pl.mouser.com/ProductDetail/OSRAM-Opto-Semiconductors/KRTBDWLM3232-T4U6-JW%2bAAAB-J3%2bR4S6-5V?qs=DPoM0jnrROWJ%2F8aeLsBacQ%3D%3D"

$diodes.Rows.Item(10).RowHeight = 31.5
$diodes.Range("D10").Font.Size = 24
$diodes.Range("D10").Font.Bold = $true
$diodes.Range("D10").WrapText = $true
$diodes.Range("D10").VerticalAlignment = -4108

$diodes.Range("C10").Select()
